$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15 data (Model_8, id = 471)
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 471
$ws.Range("C15").Value = "scaled speed of id = 471" + [char]10 + "scaled speed of id = 470"
$ws.Range("D15").Value = "lstm(50)+do(.3)" + [char]10 + "lstm/50)+do(.3)" + [char]10 + "lstm/33)"
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = "1h back" + [char]10 + "1h forward"
$ws.Range("G15").Value = "1 week"
$ws.Range("H15").Value = "Feb March April May"
$ws.Range("I15").Value = "First 7 days of June"
$ws.Range("J15").Value = 18.2
$ws.Range("K15").Value = 16.18
$ws.Range("L15").Value = 26.47
$ws.Range("M15").Value = "Adding close sensor sequence didn't change anything. Nothing seems like changing anything."

# Match formatting of the other multi-line/wrapped comment cells in the table
$ws.Range("C15").WrapText = $true
$ws.Range("D15").WrapText = $true
$ws.Range("F15").WrapText = $true
$ws.Range("G15").WrapText = $true
$ws.Range("H15").WrapText = $true
$ws.Range("I15").WrapText = $true
$ws.Range("J15").WrapText = $true
$ws.Range("K15").WrapText = $true
$ws.Range("M15").WrapText = $true

# Row height matches the other rows holding similarly long comments
$ws.Rows.Item(15).RowHeight = 33

# Move the active selection the way it ends up after entering the new row
$ws.Range("M16").Select()
